# "mise à jour industrie"
#
# - 0D!B6: was "=1/1000" (value 1E-3) -> literal value 1 (the 1/1000
#   scaling is now applied directly to the Vecteurs percentages instead).
# - Vecteurs!C2:C13: percentage-like numbers (60, 6, 200, 30, ...) are
#   rescaled to fractions (0.06, 0.006, 0.2, 0.03, ...) i.e. divided by 1000,
#   and their custom formatting is cleared back to the workbook defaults.
# - Active sheet moves from "Production_system" to "Vecteurs".

$wb = $excel.ActiveWorkbook

# --- 0D sheet -------------------------------------------------------
$ws0D = $wb.Worksheets.Item("0D")
$ws0D.Activate()
$ws0D.Range("B6").Value = 1
$ws0D.Range("B7").Select()

# --- Vecteurs sheet ---------------------------------------------------
$wsProd = $wb.Worksheets.Item("Production_system")
$wsVec  = $wb.Worksheets.Item("Vecteurs")

# Rescale the value column (percent -> fraction, i.e. /1000) for both the
# 2020 and the 2030 blocks.
$vecRows = 2..13
foreach ($r in $vecRows) {
    $cell = $wsVec.Cells.Item($r, 3)
    $old = $cell.Value2
    $cell.Value = $old / 1000
}

# Drop the bespoke formatting that had accumulated on A1:C13 and restore
# the sheet's normal look: bold boxed labels in column A (same style as
# the other sheets' header row), bold centered years in column B, and
# plain (default) everything else.
$wsVec.Range("A1:C13").ClearFormats()

$wsProd.Range("A1").Copy()
$wsVec.Range("A2:A13").PasteSpecial(-4122)

$ws0D.Range("A5").Copy()
$wsVec.Range("B2:B13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Vecteurs becomes the active sheet/selection.
$wsVec.Activate()
$wsVec.Range("G18").Select()
